# Error Calculations and Plots
# Apply the target edits to the "missing_data" worksheet:
#  1. Remove two whole data rows (RM 232 and SC 92) - everything below shifts up.
#  2. Update the "missing" mask (blank vs. numeric) on a handful of remaining cells
#     so the final grid matches the regenerated missing-data pattern.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete the two rows that were dropped from the dataset ---
# Delete from the bottom up so row numbers of earlier rows stay stable.
$ws.Rows.Item(28).Delete()   # "SC 92"
$ws.Rows.Item(26).Delete()   # "RM 232"

# --- 2. Fix up the remaining "missing value" mask cells ---
# Cells that become blank (simulated missing data -> inline empty string)
$blankCells = @("F3", "F5", "E6", "F8", "E14", "F19", "F22", "E23", "E24", "D27", "F27", "D32")
foreach ($addr in $blankCells) {
    $ws.Range($addr).ClearContents()
    $ws.Range($addr).Value = ""
}

# Cells that get a concrete numeric value (previously blank/missing)
$ws.Range("E2").Value = -7.2
$ws.Range("F4").Value = 17.97
$ws.Range("E12").Value = -5.3
$ws.Range("F15").Value = 16.2
$ws.Range("F18").Value = 18.35
$ws.Range("E20").Value = -7.2
$ws.Range("E21").Value = -8.699999999999999
$ws.Range("F23").Value = 16.48
$ws.Range("F25").Value = 16.6
$ws.Range("D26").Value = -13.8
$ws.Range("D30").Value = -13.6
$ws.Range("E31").Value = -8.1
$ws.Range("E33").Value = -10.7
